$wb = $excel.ActiveWorkbook
$dataSheet = $wb.Worksheets.Item("data")

# Update the "time_taken" timestamps on the data sheet (F2:F100)
$dataSheet.Range("F2").Value = "2021-10-05 14:33:21.192262"
$dataSheet.Range("F3").Value = "2021-10-05 14:33:21.192268"
$dataSheet.Range("F4").Value = "2021-10-05 14:33:21.192271"
$dataSheet.Range("F5").Value = "2021-10-05 14:33:21.192273"
$dataSheet.Range("F6").Value = "2021-10-05 14:33:21.192275"
$dataSheet.Range("F7").Value = "2021-10-05 14:33:21.192277"
$dataSheet.Range("F8").Value = "2021-10-05 14:33:21.192279"
$dataSheet.Range("F9").Value = "2021-10-05 14:33:21.192281"
$dataSheet.Range("F10").Value = "2021-10-05 14:33:21.192283"
$dataSheet.Range("F11").Value = "2021-10-05 14:33:21.192285"
$dataSheet.Range("F12").Value = "2021-10-05 14:33:21.192287"
$dataSheet.Range("F13").Value = "2021-10-05 14:33:21.192289"
$dataSheet.Range("F14").Value = "2021-10-05 14:33:21.192291"
$dataSheet.Range("F15").Value = "2021-10-05 14:33:21.192293"
$dataSheet.Range("F16").Value = "2021-10-05 14:33:21.192294"
$dataSheet.Range("F17").Value = "2021-10-05 14:33:21.192296"
$dataSheet.Range("F18").Value = "2021-10-05 14:33:21.192298"
$dataSheet.Range("F19").Value = "2021-10-05 14:33:21.192301"
$dataSheet.Range("F20").Value = "2021-10-05 14:33:21.192303"
$dataSheet.Range("F21").Value = "2021-10-05 14:33:21.192304"
$dataSheet.Range("F22").Value = "2021-10-05 14:33:21.192306"
$dataSheet.Range("F23").Value = "2021-10-05 14:33:21.192308"
$dataSheet.Range("F24").Value = "2021-10-05 14:33:21.192310"
$dataSheet.Range("F25").Value = "2021-10-05 14:33:21.192312"
$dataSheet.Range("F26").Value = "2021-10-05 14:33:21.192314"
$dataSheet.Range("F27").Value = "2021-10-05 14:33:21.192316"
$dataSheet.Range("F28").Value = "2021-10-05 14:33:21.192318"
$dataSheet.Range("F29").Value = "2021-10-05 14:33:21.192320"
$dataSheet.Range("F30").Value = "2021-10-05 14:33:21.192322"
$dataSheet.Range("F31").Value = "2021-10-05 14:33:21.192323"
$dataSheet.Range("F32").Value = "2021-10-05 14:33:21.192325"
$dataSheet.Range("F33").Value = "2021-10-05 14:33:21.192327"
$dataSheet.Range("F34").Value = "2021-10-05 14:33:21.192329"
$dataSheet.Range("F35").Value = "2021-10-05 14:33:21.192331"
$dataSheet.Range("F36").Value = "2021-10-05 14:33:21.192333"
$dataSheet.Range("F37").Value = "2021-10-05 14:33:21.192335"
$dataSheet.Range("F38").Value = "2021-10-05 14:33:21.192337"
$dataSheet.Range("F39").Value = "2021-10-05 14:33:21.192339"
$dataSheet.Range("F40").Value = "2021-10-05 14:33:21.192341"
$dataSheet.Range("F41").Value = "2021-10-05 14:33:21.192343"
$dataSheet.Range("F42").Value = "2021-10-05 14:33:21.192345"
$dataSheet.Range("F43").Value = "2021-10-05 14:33:21.192346"
$dataSheet.Range("F44").Value = "2021-10-05 14:33:21.192348"
$dataSheet.Range("F45").Value = "2021-10-05 14:33:21.192350"
$dataSheet.Range("F46").Value = "2021-10-05 14:33:21.192352"
$dataSheet.Range("F47").Value = "2021-10-05 14:33:21.192354"
$dataSheet.Range("F48").Value = "2021-10-05 14:33:21.192355"
$dataSheet.Range("F49").Value = "2021-10-05 14:33:21.192357"
$dataSheet.Range("F50").Value = "2021-10-05 14:33:21.192359"
$dataSheet.Range("F51").Value = "2021-10-05 14:33:21.192361"
$dataSheet.Range("F52").Value = "2021-10-05 14:33:21.192363"
$dataSheet.Range("F53").Value = "2021-10-05 14:33:21.192365"
$dataSheet.Range("F54").Value = "2021-10-05 14:33:21.192367"
$dataSheet.Range("F55").Value = "2021-10-05 14:33:21.192369"
$dataSheet.Range("F56").Value = "2021-10-05 14:33:21.192371"
$dataSheet.Range("F57").Value = "2021-10-05 14:33:21.192373"
$dataSheet.Range("F58").Value = "2021-10-05 14:33:21.192375"
$dataSheet.Range("F59").Value = "2021-10-05 14:33:21.192376"
$dataSheet.Range("F60").Value = "2021-10-05 14:33:21.192378"
$dataSheet.Range("F61").Value = "2021-10-05 14:33:21.192380"
$dataSheet.Range("F62").Value = "2021-10-05 14:33:21.192382"
$dataSheet.Range("F63").Value = "2021-10-05 14:33:21.192384"
$dataSheet.Range("F64").Value = "2021-10-05 14:33:21.192386"
$dataSheet.Range("F65").Value = "2021-10-05 14:33:21.192388"
$dataSheet.Range("F66").Value = "2021-10-05 14:33:21.192391"
$dataSheet.Range("F67").Value = "2021-10-05 14:33:21.192394"
$dataSheet.Range("F68").Value = "2021-10-05 14:33:21.192396"
$dataSheet.Range("F69").Value = "2021-10-05 14:33:21.192398"
$dataSheet.Range("F70").Value = "2021-10-05 14:33:21.192400"
$dataSheet.Range("F71").Value = "2021-10-05 14:33:21.192401"
$dataSheet.Range("F72").Value = "2021-10-05 14:33:21.192403"
$dataSheet.Range("F73").Value = "2021-10-05 14:33:21.192405"
$dataSheet.Range("F74").Value = "2021-10-05 14:33:21.192407"
$dataSheet.Range("F75").Value = "2021-10-05 14:33:21.192409"
$dataSheet.Range("F76").Value = "2021-10-05 14:33:21.192411"
$dataSheet.Range("F77").Value = "2021-10-05 14:33:21.192413"
$dataSheet.Range("F78").Value = "2021-10-05 14:33:21.192416"
$dataSheet.Range("F79").Value = "2021-10-05 14:33:21.192418"
$dataSheet.Range("F80").Value = "2021-10-05 14:33:21.192420"
$dataSheet.Range("F81").Value = "2021-10-05 14:33:21.192422"
$dataSheet.Range("F82").Value = "2021-10-05 14:33:21.192424"
$dataSheet.Range("F83").Value = "2021-10-05 14:33:21.192426"
$dataSheet.Range("F84").Value = "2021-10-05 14:33:21.192428"
$dataSheet.Range("F85").Value = "2021-10-05 14:33:21.192430"
$dataSheet.Range("F86").Value = "2021-10-05 14:33:21.192432"
$dataSheet.Range("F87").Value = "2021-10-05 14:33:21.192434"
$dataSheet.Range("F88").Value = "2021-10-05 14:33:21.192436"
$dataSheet.Range("F89").Value = "2021-10-05 14:33:21.192438"
$dataSheet.Range("F90").Value = "2021-10-05 14:33:21.192440"
$dataSheet.Range("F91").Value = "2021-10-05 14:33:21.192442"
$dataSheet.Range("F92").Value = "2021-10-05 14:33:21.192444"
$dataSheet.Range("F93").Value = "2021-10-05 14:33:21.192446"
$dataSheet.Range("F94").Value = "2021-10-05 14:33:21.192448"
$dataSheet.Range("F95").Value = "2021-10-05 14:33:21.192451"
$dataSheet.Range("F96").Value = "2021-10-05 14:33:21.192453"
$dataSheet.Range("F97").Value = "2021-10-05 14:33:21.192455"
$dataSheet.Range("F98").Value = "2021-10-05 14:33:21.192457"
$dataSheet.Range("F99").Value = "2021-10-05 14:33:21.192459"
$dataSheet.Range("F100").Value = "2021-10-05 14:33:21.192461"

# Add the new "metadata" sheet right after "data"
$ws2 = $wb.Worksheets.Add($null, $dataSheet)
$ws2.Name = "metadata"

# Header row (bold / centered / bordered -- same look as the "data" sheet header)
$ws2.Range("B1").Value = "data_name"
$ws2.Range("C1").Value = "data_id"
$ws2.Range("D1").Value = "data_version"
$ws2.Range("E1").Value = "data_version_created"
$ws2.Range("F1").Value = "panel_query_time"
$ws2.Range("G1").Value = "panel_get_request"

$dataSheet.Range("B1").Copy()
$ws2.Range("B1:G1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Data row
$ws2.Range("A2").Value = 0
$dataSheet.Range("B1").Copy()
$ws2.Range("A2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws2.Range("B2").Value = "Cancer Predisposition_Paediatric"
$ws2.Range("C2").Value = 152

# Force D2 to stay text ("0.113") instead of being coerced to a number,
# then drop back to the default cell style so no stray number format sticks around.
$ws2.Range("D2").NumberFormat = "@"
$ws2.Range("D2").Value = "0.113"
$ws2.Range("D2").Style = "Normal"

$ws2.Range("E2").Value = "2021-08-31T03:01:28.230561Z"
$ws2.Range("F2").Value = "2021-10-05 14:33:21.189525"
$ws2.Range("G2").Value = "https://panelapp.agha.umccr.org/api/v1/panels/152/?format=json"
